# Give Creatures Random Names
# The "Creature Cannabalism" task row (row 19) is removed from the Gantt
# chart; remaining rows shift up by one. The task now above it ("Creature
# Names", which becomes the new row 19) gets its schedule bar extended by
# one more week (column AQ), reflecting the extra "download list of names /
# random first name / surname" work described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the look of an existing "scheduled week" marker cell (black fill,
# date number format) before row numbers shift, so we can stamp the same
# format onto the new cell after the delete.
$ws.Range("AP18").Copy()

# Remove the whole row for "    Creature Cannabalism" (old row 19). Excel
# automatically renumbers every row below, drops the now-unused shared
# string, and shrinks the sheet dimension.
$ws.Rows("19:19").Delete()

# Stamp the black "scheduled" format onto AQ19 (one column further right
# than the previous last marker on that task's row), matching the extended
# bar for "    Creature Names" in its new position as row 19.
$ws.Range("AQ19").PasteSpecial(-4122)

# Leave the selection on the row that was just edited, the same way Excel
# leaves focus on the row after a delete/row-edit operation.
$ws.Rows("19:19").Select() | Out-Null
